# Insert a new price record at row 233 (weekly "Fruta / hortaliza" update).
# Every existing row from 233-253 shifts down by one (to 234-254); the
# new row is populated with the latest "Ciruela - Angeleno" reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("233:233").Insert()

$ws.Range("A233").Value = 10
$ws.Range("B233").Value = "Vega Modelo de Temuco"
$ws.Range("C233").Value = "La Araucanía"
$ws.Range("D233").Value = 44918
$ws.Range("E233").Value = 9
$ws.Range("F233").Value = "Fruta"
$ws.Range("G233").Value = 100103
$ws.Range("H233").Value = "Frutos de hueso (carozo)"
$ws.Range("I233").Value = 100103002
$ws.Range("J233").Value = "Ciruela"
$ws.Range("K233").Value = "Angeleno"
$ws.Range("L233").Value = "Primera"
$ws.Range("M233").Value = 80
$ws.Range("N233").Value = 12000
$ws.Range("O233").Value = 12000
$ws.Range("P233").Value = 12000
$ws.Range("Q233").Value = "$/caja 15 kilos granel"
$ws.Range("R233").Value = "Región de O'Higgins"
$ws.Range("S233").Value = 800
$ws.Range("T233").Value = 15
